$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last updated" timestamp footer (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 27 de Mayo de 2020 a las 06:35"

# --- India row (row 13) stats refresh ---
$ws.Range("B13").Value = 151876
$ws.Range("C13").Value = 1083
$ws.Range("D13").Value = 64426
$ws.Range("E13").Value = 83104
$ws.Range("G13").Value = 2
$ws.Range("H13").Value = 4346

# --- Row 89 stats refresh (Casos activos / Recuperados) ---
$ws.Range("D89").Value = 793
$ws.Range("E89").Value = 1213

# --- Kirguistan moves up in rank (new stats push it above Eslovaquia,
#     Nueva Zelanda and Eslovenia), so rows 97-100 shift down one slot
#     and row 97 gets Kirguistan's refreshed numbers ---
$ws.Range("A97").Value = "Kirguistan"
$ws.Range("B97").Value = 1520
$ws.Range("C97").Value = 52
$ws.Range("D97").Value = 1043
$ws.Range("E97").Value = 461
$ws.Range("H97").Value = 16

$ws.Range("A98").Value = "Eslovaquia"
$ws.Range("B98").Value = 1513
$ws.Range("C98").Value = 0
$ws.Range("D98").Value = 1322
$ws.Range("E98").Value = 163
$ws.Range("H98").Value = 28

$ws.Range("A99").Value = "Nueva Zelanda"
$ws.Range("B99").Value = 1504
$ws.Range("C99").Value = 0
$ws.Range("D99").Value = 1462
$ws.Range("E99").Value = 21
$ws.Range("H99").Value = 21

$ws.Range("A100").Value = "Eslovenia"
$ws.Range("B100").Value = 1469
$ws.Range("C100").Value = 0
$ws.Range("D100").Value = 1346
$ws.Range("E100").Value = 15
$ws.Range("H100").Value = 108
